$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.207906656956141
$ws.Range("C2").Value = 0.2612727038401204
$ws.Range("D2").Value = 0.03560839030920349
$ws.Range("F2").Value = 1.491134305621458
$ws.Range("G2").Value = 0.00243532404439037
$ws.Range("M2").Value = 0.889982664372809
$ws.Range("N2").Value = 1.487873516028912
$ws.Range("B3").Value = 1.086988433806766
$ws.Range("C3").Value = 0.2281047212326826
$ws.Range("D3").Value = 0.03554262974331479
$ws.Range("F3").Value = 1.4297873275641
$ws.Range("G3").Value = 0.00244100016583757
$ws.Range("M3").Value = 0.7920376295144536
$ws.Range("N3").Value = 1.491476075969715
$ws.Range("B4").Value = 1.013330230173608
$ws.Range("C4").Value = 0.2077876783929469
$ws.Range("D4").Value = 0.03552148865344762
$ws.Range("F4").Value = 1.393126634143044
$ws.Range("G4").Value = 0.002444664196820207
$ws.Range("M4").Value = 0.732455050000425
$ws.Range("N4").Value = 1.49424478498554
$ws.Range("B5").Value = 0.9834587638505923
$ws.Range("C5").Value = 0.1995196642209862
$ws.Range("D5").Value = 0.03551764733011353
$ws.Range("F5").Value = 1.37843730046238
$ws.Range("G5").Value = 0.002446202467672365
$ws.Range("M5").Value = 0.7083076036762606
$ws.Range("N5").Value = 1.495511904104035
$ws.Range("B6").Value = 0.9785073053806173
$ws.Range("C6").Value = 0.1981474293870917
$ws.Range("D6").Value = 0.03551729597500142
$ws.Range("F6").Value = 1.376013174347818
$ws.Range("G6").Value = 0.002446460628257819
$ws.Range("M6").Value = 0.7043057690154626
$ws.Range("N6").Value = 1.495730661047702
$ws.Range("B7").Value = 1.012926789943947
$ws.Range("C7").Value = 0.2076761281513768
$ws.Range("D7").Value = 0.03552141760192029
$ws.Range("F7").Value = 1.39292751944528
$ws.Range("G7").Value = 0.002444684759512592
$ws.Range("M7").Value = 0.7321288598389657
$ws.Range("N7").Value = 1.494261313107856
$ws.Range("B8").Value = 1.166090845071039
$ws.Range("C8").Value = 0.2498258888563498
$ws.Range("D8").Value = 0.03558168521084326
$ws.Range("F8").Value = 1.469771158553655
$ws.Range("G8").Value = 0.002437244149819325
$ws.Range("M8").Value = 0.8560915894332624
$ws.Range("N8").Value = 1.488999437483457
$ws.Range("B9").Value = 1.471226435440769
$ws.Range("C9").Value = 0.3329038136561735
$ws.Range("D9").Value = 0.03585541601623277
$ws.Range("F9").Value = 1.628590148478366
$ws.Range("G9").Value = 0.002424064575472778
$ws.Range("M9").Value = 1.103922011171832
$ws.Range("N9").Value = 1.483150767534752
$ws.Range("B10").Value = 1.698530926909143
$ws.Range("C10").Value = 0.3942598878412014
$ws.Range("D10").Value = 0.03615548230681043
$ws.Range("F10").Value = 1.75044133117072
$ws.Range("G10").Value = 0.002415231021234734
$ws.Range("M10").Value = 1.289371578000257
$ws.Range("N10").Value = 1.481652699949436
$ws.Range("B11").Value = 1.802662232735429
$ws.Range("C11").Value = 0.4222553460726317
$ws.Range("D11").Value = 0.03631439266431613
$ws.Range("F11").Value = 1.807043689552103
$ws.Range("G11").Value = 0.002411394509811818
$ws.Range("M11").Value = 1.37457661968169
$ws.Range("N11").Value = 1.481595273199318
$ws.Range("B12").Value = 1.8422023171305
$ws.Range("C12").Value = 0.4328696295687564
$ws.Range("D12").Value = 0.03637786600428683
$ws.Range("F12").Value = 1.82864964613205
$ws.Range("G12").Value = 0.002409967700758379
$ws.Range("M12").Value = 1.406971522056864
$ws.Range("N12").Value = 1.481664617536126
$ws.Range("B13").Value = 1.833681809733321
$ws.Range("C13").Value = 0.4305830606359109
$ws.Range("D13").Value = 0.03636404801827808
$ws.Range("F13").Value = 1.823988718115288
$ws.Range("G13").Value = 0.00241027383597316
$ws.Range("M13").Value = 1.399988800403293
$ws.Range("N13").Value = 1.481645609921003
$ws.Range("B14").Value = 1.80591304135919
$ws.Range("C14").Value = 0.423128323626031
$ws.Range("D14").Value = 0.03631954813317861
$ws.Range("F14").Value = 1.808817760399705
$ws.Range("G14").Value = 0.002411276605352634
$ws.Range("M14").Value = 1.377239121608
$ws.Range("N14").Value = 1.481599144782763
$ws.Range("B15").Value = 1.788917998450927
$ws.Range("C15").Value = 0.4185638031750614
$ws.Range("D15").Value = 0.03629272234580583
$ws.Range("F15").Value = 1.799547595194895
$ws.Range("G15").Value = 0.002411894210728947
$ws.Range("M15").Value = 1.36332141683171
$ws.Range("N15").Value = 1.481582586751813
$ws.Range("B16").Value = 1.691740627362265
$ws.Range("C16").Value = 0.3924320905055652
$ws.Range("D16").Value = 0.03614555444718803
$ws.Range("F16").Value = 1.746766098920489
$ws.Range("G16").Value = 0.002415485392818167
$ws.Range("M16").Value = 1.28382090325016
$ws.Range("N16").Value = 1.481669125562732
$ws.Range("B17").Value = 1.632314333118813
$ws.Range("C17").Value = 0.3764233451595942
$ws.Range("D17").Value = 0.03606106408029319
$ws.Range("F17").Value = 1.71468884072911
$ws.Range("G17").Value = 0.002417734942306341
$ws.Range("M17").Value = 1.23527173017338
$ws.Range("N17").Value = 1.481883004373117
$ws.Range("B18").Value = 1.598202370884508
$ws.Range("C18").Value = 0.3672234159230925
$ws.Range("D18").Value = 0.0360145720901599
$ws.Range("F18").Value = 1.69634892461292
$ws.Range("G18").Value = 0.002419045956344771
$ws.Range("M18").Value = 1.20742637238169
$ws.Range("N18").Value = 1.48206467830515
$ws.Range("B19").Value = 1.586664317975647
$ws.Range("C19").Value = 0.3641098040418456
$ws.Range("D19").Value = 0.03599918991064044
$ws.Range("F19").Value = 1.690158147153937
$ws.Range("G19").Value = 0.00241949279089428
$ws.Range("M19").Value = 1.198011702931964
$ws.Range("N19").Value = 1.482136223410293
$ws.Range("B20").Value = 1.638633256074797
$ws.Range("C20").Value = 0.3781266817089772
$ws.Range("D20").Value = 0.03606983986316692
$ws.Range("F20").Value = 1.718092102971696
$ws.Range("G20").Value = 0.002417493701765125
$ws.Range("M20").Value = 1.240431652782064
$ws.Range("N20").Value = 1.481854156159358
$ws.Range("B21").Value = 1.814066451840347
$ws.Range("C21").Value = 0.4253176002746386
$ws.Range("D21").Value = 0.03633252873116533
$ws.Range("F21").Value = 1.813269146554603
$ws.Range("G21").Value = 0.002410981363664721
$ws.Range("M21").Value = 1.383917668522628
$ws.Range("N21").Value = 1.481610309652666
$ws.Range("B22").Value = 1.929352047894326
$ws.Range("C22").Value = 0.4562358587294284
$ws.Range("D22").Value = 0.03652346697434439
$ws.Range("F22").Value = 1.876476264070448
$ws.Range("G22").Value = 0.002406876612615927
$ws.Range("M22").Value = 1.478453549018141
$ws.Range("N22").Value = 1.481982531468944
$ws.Range("B23").Value = 1.867763335255688
$ws.Range("C23").Value = 0.4397269245177426
$ws.Range("D23").Value = 0.0364197725291362
$ws.Range("F23").Value = 1.842648487084659
$ws.Range("G23").Value = 0.002409053593287399
$ws.Range("M23").Value = 1.427925652409428
$ws.Range("N23").Value = 1.481734776344041
$ws.Range("B24").Value = 1.635776306126047
$ws.Range("C24").Value = 0.3773565917610426
$ws.Range("D24").Value = 0.03606586585068072
$ws.Range("F24").Value = 1.71655317151081
$ws.Range("G24").Value = 0.00241760271147501
$ws.Range("M24").Value = 1.238098646175729
$ws.Range("N24").Value = 1.481867015663838
$ws.Range("B25").Value = 1.388146166320382
$ws.Range("C25").Value = 0.3103775442822041
$ws.Range("D25").Value = 0.03576429005553194
$ws.Range("F25").Value = 1.584733506319594
$ws.Range("G25").Value = 0.002427480030103785
$ws.Range("M25").Value = 1.036321109343234
$ws.Range("N25").Value = 1.484246881604037
